$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

$ws.Cells.Item($row, 1).Value = 61
$ws.Cells.Item($row, 2).Value = "armenia"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45225.70833333334
$ws.Cells.Item($row, 6).Value = "Shirak Gyumri"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "Noah"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 4.51
$ws.Cells.Item($row, 11).Value = "25/10/2023 04:14"
$ws.Cells.Item($row, 12).Value = 5.9
$ws.Cells.Item($row, 13).Value = "26/10/2023 16:59"
$ws.Cells.Item($row, 14).Value = 3.67
$ws.Cells.Item($row, 15).Value = "25/10/2023 04:14"
$ws.Cells.Item($row, 16).Value = 3.67
$ws.Cells.Item($row, 17).Value = "26/10/2023 16:59"
$ws.Cells.Item($row, 18).Value = 1.67
$ws.Cells.Item($row, 19).Value = "25/10/2023 04:14"
$ws.Cells.Item($row, 20).Value = 1.63
$ws.Cells.Item($row, 21).Value = "26/10/2023 16:59"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/armenia/premier-league/shirak-gyumri-noah/EFOOAPD8/"

# Copy the formatting from the previous row (61) so the new row matches
# the style of the rest of the table (bold/bordered/centered column A,
# datetime-formatted column E).
$ws.Range("A61").Copy() | Out-Null
$ws.Range("A" + $row).PasteSpecial(-4122) | Out-Null
$ws.Range("E61").Copy() | Out-Null
$ws.Range("E" + $row).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
